$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.004.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.63%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.580.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.08%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +9.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '570.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("E7").Value = '  -0.83%  '

$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("E9").Value = '  +0.32%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '61.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +11.60%  '

$ws.Range("E11").Value = '  -1.79%  '

$ws.Range("E12").Value = '  +4.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.24'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.148.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.579.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("E16").Value = '  +0.63%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.857.03'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.63%  '

$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("E21").Value = '  +1.65%  '

$ws.Range("E22").Value = '  -0.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +11.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.42%  '

$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.72%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.25%  '

$ws.Range("E27").Value = '  -0.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.49'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '676.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.47%  '

$ws.Range("E32").Value = '  -0.86%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.112'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.58%  '

$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.07'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.82%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '40.93'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.409'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.48%  '

$ws.Range("E37").Value = '  -0.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.20'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.86%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0748'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.170.26'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.85%  '

$ws.Range("E41").Value = '  +0.29%  '

$ws.Range("E42").Value = '  -0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +13.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.82%  '

$ws.Range("E46").Value = '  -1.06%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.22%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.130'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.21%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '138.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.56%  '

$ws.Range("E51").Value = '  -1.90%  '
